# Re-curate the metadata for the "municipio-nombre" (col B) and
# "intervalo-renta" (col C) dimensions, per the new curated dimension set.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: property URIs
$ws.Range("B2").Value = "sdmx-dimension:refArea"          # was iaest-measure:municipio-nombre
$ws.Range("C2").Value = "iaest-measure:intervalo-renta"   # was iaest-dimension:intervalo-renta

# Row 3: dim/medida role swaps between the two columns
$ws.Range("B3").Value = "dim"      # was medida
$ws.Range("C3").Value = "medida"   # was dim

# Row 4: rdfs/xsd type swaps between the two columns
$ws.Range("B4").Value = "URI-Municipio"   # was xsd:int
$ws.Range("C4").Value = "xsd:int"         # was skos:Concept

# Row 5: intervalo-renta no longer has an associated mapping workbook
$ws.Range("C5").Clear()
